# "Generate Report for Handback"
#
# This script updates the handback-status report:
#  - The overall handback status text changes from
#    "Handed back: in sync with en-US" to
#    "Handed back: not in sync with en-US" (shown on the Overview sheet
#    as well as the Status column of the zh-cn and de-de detail sheets).
#  - A new "Correspond Handback DateTime" is recorded for the second
#    file (7990fd38-cd2c-40ec-a0bf-423934b0b2cd) on both the zh-cn and
#    de-de sheets, reflecting a fresh handback run.
#  - The Status column (and the Overview sheet's language columns) are
#    widened slightly to fit the longer status text.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: not in sync with en-US"

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# Update the status text everywhere it appears.
$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus
$ws1.Range("E3").Value = $newStatus
$ws1.Range("F3").Value = $newStatus

$ws2.Range("C2").Value = $newStatus
$ws2.Range("C3").Value = $newStatus

$ws3.Range("C2").Value = $newStatus
$ws3.Range("C3").Value = $newStatus

# Record the new handback datetimes for the second file on each
# language sheet ("Correspond Handback DateTime" column K).
$ws2.Range("K3").Value = "2016-11-15 17:45:54"
$ws3.Range("K3").Value = "2016-11-15 17:46:13"

# Widen the columns that now hold the longer status text so the
# report still displays cleanly (mirrors Excel auto-fit behavior).
$ws1.Range("E1").ColumnWidth = 32.66
$ws1.Range("F1").ColumnWidth = 32.66
$ws2.Range("C1").ColumnWidth = 32.66
$ws3.Range("C1").ColumnWidth = 32.66
